# Applies the "Constant trace and other changes" edit to RS.xlsx
#
# Summary of changes:
#  - Row 2 (R1,i), Row 3 (S1,i), Row 4 (T1,i), Row 7 (T2,i), Row 10 (T3,i):
#    previously-empty trailing cells across B:V are filled in with constant
#    (trace) values, mostly 0, with a handful of non-zero entries.
#  - Row 6 (S2,i) and Row 9 (S3,i): the "-" placeholder text in columns B and V
#    is replaced with actual numeric values (so the now-unused "-" shared
#    string is dropped automatically on save).
#  - Cell E3 gets a dedicated "0.00" number format (distinct from the sheet's
#    usual "0.000" format).
#  - The worksheet view selection becomes the A1:V10 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : R1,i --------------------------------------------------------
$ws.Cells.Item(2, 2).Value  = 3.3661555587195302   # B2
$ws.Cells.Item(2, 3).Value  = -11.2586386586866    # C2
$ws.Cells.Item(2, 4).Value  = 15.533223141117199   # D2
$ws.Cells.Item(2, 5).Value  = -10.6631010062975    # E2
$ws.Cells.Item(2, 6).Value  = 3.0227262326444699   # F2
for ($col = 7; $col -le 22; $col++) {
    $ws.Cells.Item(2, $col).Value = 0               # G2:V2
}

# ---- Row 3 : S1,i ---------------------------------------------------------
# B3 already equals 1; fill in the rest.
$ws.Cells.Item(3, 3).Value = -1                     # C3
$ws.Cells.Item(3, 4).Value = 0                       # D3
$ws.Cells.Item(3, 5).Value = 0                       # E3 (gets special format below)
for ($col = 6; $col -le 22; $col++) {
    $ws.Cells.Item(3, $col).Value = 0                 # F3:V3
}
# E3 uses its own "0.00" number format (distinct from the default "0.000").
$ws.Cells.Item(3, 5).NumberFormat = "0.00"

# ---- Row 4 : T1,i ----------------------------------------------------------
$ws.Cells.Item(4, 2).Value = 3.3661555587195302     # B4
$ws.Cells.Item(4, 3).Value = -11.2586386586866      # C4
$ws.Cells.Item(4, 4).Value = 15.533223141117199     # D4
$ws.Cells.Item(4, 5).Value = -10.6631010062975      # E4
$ws.Cells.Item(4, 6).Value = 3.0227262326444699     # F4
for ($col = 7; $col -le 22; $col++) {
    $ws.Cells.Item(4, $col).Value = 0                 # G4:V4
}

# ---- Row 6 : S2,i -----------------------------------------------------------
# B6 and V6 previously held the text placeholder "-"; now numeric values.
$ws.Cells.Item(6, 2).Value = 1                       # B6
$ws.Cells.Item(6, 22).Value = 0                      # V6

# ---- Row 7 : T2,i ------------------------------------------------------------
# B7 already equals 0.198; fill in the rest with zeros.
for ($col = 3; $col -le 22; $col++) {
    $ws.Cells.Item(7, $col).Value = 0                 # C7:V7
}

# ---- Row 9 : S3,i ------------------------------------------------------------
# B9 and V9 previously held the text placeholder "-"; now numeric values.
$ws.Cells.Item(9, 2).Value = 0                       # B9
$ws.Cells.Item(9, 22).Value = 0                      # V9

# ---- Row 10 : T3,i -----------------------------------------------------------
# B10 already equals 0.137; fill in the rest with zeros. These cells did not
# exist before, so they need the same "0.000" number format as the rest of
# the data rows (matching style index "1" used throughout columns B:V).
for ($col = 3; $col -le 22; $col++) {
    $ws.Cells.Item(10, $col).Value = 0                # C10:V10
    $ws.Cells.Item(10, $col).NumberFormat = "0.000"
}

# ---- View selection -----------------------------------------------------------
$ws.Range("A1:V10").Select()
